$d = $word.ActiveDocument

$d.Content.Find.Execute("945÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "104÷4=", 2) | Out-Null
$d.Content.Find.Execute("531÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "171÷4=", 2) | Out-Null
$d.Content.Find.Execute("725÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "318÷8=", 2) | Out-Null
$d.Content.Find.Execute("104÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "149÷5=", 2) | Out-Null
$d.Content.Find.Execute("536÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "773÷5=", 2) | Out-Null
$d.Content.Find.Execute("160÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "296÷9=", 2) | Out-Null
$d.Content.Find.Execute("904÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "233÷8=", 2) | Out-Null
$d.Content.Find.Execute("909÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "375÷7=", 2) | Out-Null
$d.Content.Find.Execute("575÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "335÷3=", 2) | Out-Null
$d.Content.Find.Execute("768÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "752÷7=", 2) | Out-Null
$d.Content.Find.Execute("705÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "573÷2=", 2) | Out-Null
$d.Content.Find.Execute("624÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "662÷3=", 2) | Out-Null
$d.Content.Find.Execute("570÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "769÷5=", 2) | Out-Null
$d.Content.Find.Execute("620÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "468÷7=", 2) | Out-Null
$d.Content.Find.Execute("878÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "568÷2=", 2) | Out-Null
$d.Content.Find.Execute("433÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "665÷7=", 2) | Out-Null
$d.Content.Find.Execute("701÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "453÷3=", 2) | Out-Null
$d.Content.Find.Execute("941÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "599÷4=", 2) | Out-Null
$d.Content.Find.Execute("311÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "870÷3=", 2) | Out-Null
$d.Content.Find.Execute("554÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "422÷4=", 2) | Out-Null
$d.Content.Find.Execute("481÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "606÷3=", 2) | Out-Null
$d.Content.Find.Execute("948÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "688÷3=", 2) | Out-Null
$d.Content.Find.Execute("418÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "158÷9=", 2) | Out-Null
$d.Content.Find.Execute("863÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "162÷8=", 2) | Out-Null
$d.Content.Find.Execute("282÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "665÷9=", 2) | Out-Null
